$wb = $excel.ActiveWorkbook

# --- ModuleVariables sheet: update Heater-Shaker labware API name ---
$wsModule = $wb.Worksheets.Item("ModuleVariables")
$wsModule.Range("B8").Value = "3dprinted_opentrons_shaker_1.5mleppendorf"

# --- GeneralVariables sheet: update labware API names ---
$wsGeneral = $wb.Worksheets.Item("GeneralVariables")
$wsGeneral.Range("B2").Value = "nest_96_wellplate_200ul_flat"
$wsGeneral.Range("B4").Value = "biorad_96_wellplate_200ul_pcr"

# --- Restore on-screen selections: cursor left on GeneralVariables!B4 after
# the edit, but the workbook was last saved with ModuleVariables as the
# active tab (selection B8), so re-select that sheet/cell last. ---
$wsGeneral.Range("B4").Select()
$wsModule.Range("B8").Select()
